# Insert a new data row at row 302 (weekly "Ajo" / garlic price record),
# shifting the existing rows 302-360 down to 303-361.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(302).Insert()

$ws.Cells.Item(302, 1).Value  = 9
$ws.Cells.Item(302, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(302, 3).Value  = "Metropolitana"
$ws.Cells.Item(302, 4).Value  = (Get-Date -Year 2023 -Month 10 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(302, 5).Value  = 13
$ws.Cells.Item(302, 6).Value  = 100112003
$ws.Cells.Item(302, 7).Value  = "Ajo"
$ws.Cells.Item(302, 8).Value  = "Chino"
$ws.Cells.Item(302, 9).Value  = "Primera"
$ws.Cells.Item(302, 10).Value = 340
$ws.Cells.Item(302, 11).Value = 18000
$ws.Cells.Item(302, 12).Value = 19000
$ws.Cells.Item(302, 13).Value = 18500
$ws.Cells.Item(302, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(302, 15).Value = "China"
$ws.Cells.Item(302, 16).Value = 1850
$ws.Cells.Item(302, 17).Value = 10
$ws.Cells.Item(302, 18).Value = "Hortaliza"
